$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the header labels for columns B, C, D
$ws.Range("B1").Value = "target_linear"
$ws.Range("C1").Value = "target_ambitious"
$ws.Range("D1").Value = "target_central"

# Delete the now-unused columns E and F (mild_s / agressive_s data)
$ws.Range("E1:F32").Delete(-4161)

$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.Refresh()

$ws.Range("F19").Select()

$wb.Save()
